{"js": "// The paragraph about presentations gets two insertions:\n//  1. \"internacional \" is inserted right after \"... \u00faltimo congresso \"\n//     (so \"congresso de Sistemas Complexos\" becomes\n//      \"congresso internacional de Sistemas Complexos\").\n//  2. \" para a Secretaria-Geral da Presid\u00eancia da Rep\u00fablica\" is inserted\n//     right after \"...participativos brasileiros\" (before the period).\nconst body = context.document.body;\n\n// 1) \"...no \u00faltimo congresso de Sistemas Complexos...\"\nconst congressoResults = body.search(\"\u00faltimo congresso \", {\n  matchCase: false,\n  matchWholeWord: false\n});\ncongressoResults.load(\"items\");\nawait context.sync();\n\nif (congressoResults.items.length > 0) {\n  const congressoRange = congressoResults.items[0];\n  const insertAfterCongresso = congressoRange.getRange(\"End\");\n  insertAfterCongresso.insertText(\"internacional \", \"Before\");\n  await context.sync();\n}\n\n// 2) \"...participativos brasileiros. Todas estas atividades...\"\nconst brasileirosResults = body.search(\"participativos brasileiros\", {\n  matchCase: false,\n  matchWholeWord: false\n});\nbrasileirosResults.load(\"items\");\nawait context.sync();\n\nif (brasileirosResults.items.length > 0) {\n  const brasileirosRange = brasileirosResults.items[0];\n  const insertAfterBrasileiros = brasileirosRange.getRange(\"End\");\n  insertAfterBrasileiros.insertText(\n    \" para a Secretaria-Geral da Presid\u00eancia da Rep\u00fablica\",\n    \"Before\"\n  );\n  await context.sync();\n}\n", "ps1": "# The paragraph about presentations gets two insertions:\n#  1. \"internacional \" is inserted right after \"... \u00faltimo congresso \"\n#     (so \"congresso de Sistemas Complexos\" becomes\n#      \"congresso internacional de Sistemas Complexos\").\n#  2. \" para a Secretaria-Geral da Presid\u00eancia da Rep\u00fablica\" is inserted\n#     right after \"...participativos brasileiros\" (before the period).\n$d = $word.ActiveDocument\n\n# 1) \"...no \u00faltimo congresso de Sistemas Complexos...\"\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Execute(\n    \"congresso de Sistemas Complexos\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"congresso internacional de Sistemas Complexos\",\n    1\n) | Out-Null\n\n# 2) \"...participativos brasileiros. Todas estas atividades...\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Execute(\n    \"participativos brasileiros.\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"participativos brasileiros para a Secretaria-Geral da Presid\u00eancia da Rep\u00fablica.\",\n    1\n) | Out-Null\n"}
